# Update column A values (rows 2-86) on the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    28.30556419040565,
    7.676473805426895,
    7.944418227584919,
    3.415618938336195,
    3.116020622501082,
    2.994456602645386,
    5.88973572244214,
    1.534354667387845,
    10.12321440681649,
    0.5900651275773328,
    6.028896973072136,
    3.609487352116929,
    4.609890264698635,
    0.8522022127903881,
    1.987579744537641,
    1.257736972979956,
    5.216982649957799,
    3.788602647592455,
    0.6387614721413115,
    3.414959208046184,
    2.038917012410053,
    0.4854494366466895,
    2.808907884265096,
    1.293759942376454,
    5.710266940505818,
    5.322634715292821,
    1.271915641835164,
    4.956427830928106,
    3.458700048572354,
    6.300969794730321,
    3.947004899624091,
    4.448180521353947,
    3.040576669824759,
    1.1691830588328,
    4.192214256073981,
    3.828851996518921,
    6.237926996838752,
    5.348889142888623,
    6.490047176699989,
    6.966328863329551,
    8.692575250883891,
    1.65804762674108,
    4.869295015244063,
    0.1385215065000978,
    2.560196570020395,
    3.960542001562146,
    1.587850634691875,
    4.497191368304243,
    2.687779156573725,
    3.35918414536161,
    2.227878816095682,
    0.1157204711687996,
    2.86152046210637,
    2.700725619718668,
    3.261091502558912,
    0.04716473961630641,
    2.259787257745074,
    1.051048077765387,
    0.04896573817785566,
    1.246089703777358,
    3.732817596500496,
    7.065275994188823,
    6.788761644023396,
    0.9023699298346912,
    7.210383439699797,
    2.931111611324496,
    5.354712804239426,
    7.50997915616145,
    3.300823105730927,
    0.6462444035926467,
    4.247627865884226,
    2.414099424793022,
    6.812554763277262,
    0.8685221301803381,
    1.24962359301594,
    1.351711336541911,
    2.646163582967546,
    2.122528569092339,
    0.2238756736561527,
    3.240177278725497,
    0.3967112054478719,
    1.883477187976155,
    3.332656956713947,
    0.4007973560852918,
    2.987432433173382
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
